$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - High Dividend (Europe) / iShares
$ws.Range("D15").Value = "IE00BYYHSM20"
$ws.Range("A15").Value = "High Dividend (Europe)"
$ws.Range("E15").Value = "iShares "
$ws.Range("B15").Value = 713283
$ws.Range("C15").Value = "Europe"
$ws.Range("G15").Value = 0.28

# Row 16 - Value (EM) / iShares
$ws.Range("D16").Value = "IE00BG0SKF03"
$ws.Range("A16").Value = "Value (EM)"
$ws.Range("C16").Value = "EM"
$ws.Range("B16").Value = 719637
$ws.Range("E16").Value = "iShares"
$ws.Range("G16").Value = 0.4

# Row 17 - Momentum (EM)
$ws.Range("A17").Value = "Momentum (EM)"
$ws.Range("B17").Value = 703757
$ws.Range("C17").Value = "EM"

# Row 18 - Small-Cap (Value) (EM)
$ws.Range("A18").Value = "Small-Cap (Value) (EM)"
$ws.Range("B18").Value = 702239
$ws.Range("C18").Value = "EM"

# Restore the active selection as left by the author
$null = $ws.Range("G17").Select()
